$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily-score rows appended below the existing data (rows 71-76).
# Column A holds the date as literal text (matches the existing rows'
# inlineStr storage), so force a text number format before assigning.
$newRows = @(
    @{ Row = 71; Date = "2025-02-24"; Typ = "sleep";           C = $false; D = $false },
    @{ Row = 72; Date = "2025-02-24"; Typ = "activity";        C = $true;  D = $false },
    @{ Row = 73; Date = "2025-02-24"; Typ = "weekly_activity"; C = $false; D = $false },
    @{ Row = 74; Date = "2025-02-25"; Typ = "sleep";           C = $true;  D = $false },
    @{ Row = 75; Date = "2025-02-25"; Typ = "activity";        C = $false; D = $false },
    @{ Row = 76; Date = "2025-02-25"; Typ = "weekly_activity"; C = $false; D = $false }
)

foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").NumberFormat = "@"
    $ws.Range("A$($r.Row)").Value = $r.Date
    $ws.Range("B$($r.Row)").Value = $r.Typ
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
}
